$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.203.65'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.661.01'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.26'
$ws.Range("E5").Value = '  +1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5226'
$ws.Range("E6").Value = '  -1.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2667'
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06314'
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.02'
$ws.Range("E10").Value = '  -2.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07709'
$ws.Range("E11").Value = '  -1.04%  '
$ws.Range("D12").Value = '1.675.85'
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.425'
$ws.Range("E13").Value = '  -1.32%  '
$ws.Range("D14").Value = '1.889.40'
$ws.Range("E14").Value = '  -0.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5464'
$ws.Range("E15").Value = '  -1.53%  '
$ws.Range("D16").Value = '0.0₅8201'
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.80'
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").Value = '26.244.81'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.005'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.657'
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.94'
$ws.Range("E21").Value = '  -1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.13'
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("E23").Value = '  -3.79%  '
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.08'
$ws.Range("E25").Value = '  -3.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1236'
$ws.Range("E26").Value = '  -2.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.228'
$ws.Range("E27").Value = '  -2.57%  '
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.415'
$ws.Range("E29").Value = '  -1.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05999'
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.281'
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.713'
$ws.Range("E32").Value = '  +2.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.314'
$ws.Range("E33").Value = '  -3.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.631'
$ws.Range("E34").Value = '  -3.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9790'
$ws.Range("E35").Value = '  -2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.784'
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -0.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5873'
$ws.Range("E38").Value = '  +3.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01592'
$ws.Range("E39").Value = '  -2.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.947'
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8603'
$ws.Range("E41").Value = '  +0.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = '1.033.16'
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("D45").Value = '1.803.92'
$ws.Range("E45").Value = '  -1.39%  '
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '57.14'
$ws.Range("E47").Value = '  +0.25%  '
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.098'
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("E50").Value = '  -0.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.469'
$ws.Range("E51").Value = '  +0.95%  '
